$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers (a single decimal point) must be
# explicitly forced to Text first, otherwise Excel silently reinterprets the
# assigned string as a number (e.g. "23.30" -> 23.3), which would not match the
# original inline-string cell content.
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '29.095.93'
$ws.Range('E2').Value = '  -2.37%  '
$ws.Range('D3').Value = '1.849.60'
$ws.Range('E3').Value = '  -1.23%  '
$ws.Range('D4').Value = '0.9981'
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = '0.6942'
$ws.Range('E5').Value = '  -5.20%  '
$ws.Range('E6').Value = '  -1.37%  '
$ws.Range('D7').Value = '0.9991'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = '0.07768'
$ws.Range('E8').Value = '  +8.98%  '
$ws.Range('D9').Value = '0.3042'
$ws.Range('E9').Value = '  -3.08%  '
$ws.Range('D10').Value = '23.30'
$ws.Range('D11').Value = '0.08109'
$ws.Range('E11').Value = '  -0.71%  '
$ws.Range('D12').Value = '1.861.20'
$ws.Range('E12').Value = '  -2.89%  '
$ws.Range('D13').Value = '0.7255'
$ws.Range('E13').Value = '  -2.32%  '
$ws.Range('D14').Value = '5.211'
$ws.Range('E14').Value = '  -2.50%  '
$ws.Range('D15').Value = '89.06'
$ws.Range('E15').Value = '  -3.66%  '
$ws.Range('D16').Value = '29.104.88'
$ws.Range('E16').Value = '  -2.37%  '
$ws.Range('D17').Value = '5.743'
$ws.Range('E17').Value = '  -4.47%  '
$ws.Range('D18').Value = '0.000007822'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').Value = '13.18'
$ws.Range('E19').Value = '  -1.64%  '
$ws.Range('D20').Value = '236.09'
$ws.Range('E20').Value = '  -4.91%  '
$ws.Range('D21').Value = '0.9996'
$ws.Range('E21').Value = '  -0.09%  '
$ws.Range('D22').Value = '2.098.38'
$ws.Range('E22').Value = '  -1.74%  '
$ws.Range('D23').Value = '0.9996'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').Value = '7.595'
$ws.Range('E24').Value = '  -2.18%  '
$ws.Range('D25').Value = '8.979'
$ws.Range('E25').Value = '  -2.43%  '
$ws.Range('D26').Value = '161.34'
$ws.Range('E26').Value = '  -1.55%  '
$ws.Range('D27').Value = '0.1433'
$ws.Range('E27').Value = '  -7.15%  '
$ws.Range('E28').Value = '  -2.59%  '
$ws.Range('D29').Value = '1.977'
$ws.Range('E29').Value = '  -2.10%  '
$ws.Range('D30').Value = '1.402'
$ws.Range('E30').Value = '  -3.09%  '
$ws.Range('D31').Value = '4.488'
$ws.Range('E31').Value = '  -0.81%  '
$ws.Range('D32').Value = '1.485'
$ws.Range('E32').Value = '  -2.35%  '
$ws.Range('D33').Value = '4.007'
$ws.Range('E33').Value = '  -4.28%  '
$ws.Range('E34').Value = '  -1.35%  '
$ws.Range('D35').Value = '1.181'
$ws.Range('E35').Value = '  -4.15%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').Value = '0.7031'
$ws.Range('E36').Value = '  -5.17%  '
$ws.Range('B37').Value = 'Frax'
$ws.Range('C37').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D37').Value = '1.025'
$ws.Range('E37').Value = '  +2.36%  '
$ws.Range('D38').Value = '2.648'
$ws.Range('E38').Value = '  -1.83%  '
$ws.Range('D39').Value = '0.01854'
$ws.Range('E39').Value = '  -4.18%  '
$ws.Range('D40').Value = '2.668'
$ws.Range('E40').Value = '  -2.37%  '
$ws.Range('D41').Value = '0.9128'
$ws.Range('E41').Value = '  +5.23%  '
$ws.Range('D42').Value = '1.089.75'
$ws.Range('E42').Value = '  +4.39%  '
$ws.Range('D43').Value = '6.008'
$ws.Range('E43').Value = '  +0.59%  '
$ws.Range('D44').Value = '0.4265'
$ws.Range('E44').Value = '  -4.41%  '
$ws.Range('D45').Value = '70.67'
$ws.Range('E45').Value = '  -0.97%  '
$ws.Range('D46').Value = '0.9989'
$ws.Range('E46').Value = '  -0.10%  '
$ws.Range('D47').Value = '102.90'
$ws.Range('E47').Value = '  -0.97%  '
$ws.Range('D48').Value = '1.768'
$ws.Range('E48').Value = '  -2.75%  '
$ws.Range('D49').Value = '1.996.55'
$ws.Range('E49').Value = '  -1.44%  '
$ws.Range('D50').Value = '9.154'
$ws.Range('E50').Value = '  -3.73%  '
$ws.Range('D51').Value = '6.974'
$ws.Range('E51').Value = '  -6.21%  '

# Restore the default (unstyled) cell style so only the values differ from the
# original workbook - the source data carries no explicit style on these cells.
$ws.Range('D4').Style = "Normal"
$ws.Range('D5').Style = "Normal"
$ws.Range('D7').Style = "Normal"
$ws.Range('D8').Style = "Normal"
$ws.Range('D9').Style = "Normal"
$ws.Range('D10').Style = "Normal"
$ws.Range('D11').Style = "Normal"
$ws.Range('D13').Style = "Normal"
$ws.Range('D14').Style = "Normal"
$ws.Range('D15').Style = "Normal"
$ws.Range('D17').Style = "Normal"
$ws.Range('D18').Style = "Normal"
$ws.Range('D19').Style = "Normal"
$ws.Range('D20').Style = "Normal"
$ws.Range('D21').Style = "Normal"
$ws.Range('D23').Style = "Normal"
$ws.Range('D24').Style = "Normal"
$ws.Range('D25').Style = "Normal"
$ws.Range('D26').Style = "Normal"
$ws.Range('D27').Style = "Normal"
$ws.Range('D29').Style = "Normal"
$ws.Range('D30').Style = "Normal"
$ws.Range('D31').Style = "Normal"
$ws.Range('D32').Style = "Normal"
$ws.Range('D33').Style = "Normal"
$ws.Range('D35').Style = "Normal"
$ws.Range('D36').Style = "Normal"
$ws.Range('D37').Style = "Normal"
$ws.Range('D38').Style = "Normal"
$ws.Range('D39').Style = "Normal"
$ws.Range('D40').Style = "Normal"
$ws.Range('D41').Style = "Normal"
$ws.Range('D43').Style = "Normal"
$ws.Range('D44').Style = "Normal"
$ws.Range('D45').Style = "Normal"
$ws.Range('D46').Style = "Normal"
$ws.Range('D47').Style = "Normal"
$ws.Range('D48').Style = "Normal"
$ws.Range('D50').Style = "Normal"
$ws.Range('D51').Style = "Normal"
